$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 621, shifting the existing 2026/12/29.. data down by one.
$ws.Rows("621:621").Insert()

# Write the new row's values. Column A holds a date-like text value
# ("2026/01/14"); force it to be stored as text (matching the existing
# inlineStr cells in this column) instead of letting Excel auto-convert it
# to a date serial, then clear the temporary number-format override so the
# cell is left with the workbook's default style (same as its neighbours).
$ws.Range("A621").NumberFormat = "@"
$ws.Range("A621").Value = "2026/01/14"
$ws.Range("A621").ClearFormats()

$ws.Range("B621").Value = "水"
$ws.Range("C621").Value = 11
$ws.Range("D621").Value = 30
